$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 92: values mirror the existing sheet's convention of storing
# everything (including numeric-looking values like "1111") as text.
# Force text interpretation via NumberFormat "@" so "1111" isn't coerced
# to a number, then restore the default "Normal" style so no stray
# formatting is left behind on the new cells.
$row = 92

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "1111"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "Incompleto"
$ws.Range("C$row").Value = "PS3"
$ws.Range("D$row").Value = "Zerar"
